$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Cells that change type (number <-> text): copy format then value from a sibling cell, then overwrite ---
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4163)

$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4163)

$ws.Range("C16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 3

$ws.Range("D16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1

$ws.Range("E16").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 200

$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4163)

# --- Plain value-only updates (style/type unchanged) ---
$ws.Range("L15").Value = 41.666666666666
$ws.Range("N15").Value = -22.727272727272
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -38.888888888888
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = 13.274336283185
$ws.Range("L16").Value = 18.518518518518
$ws.Range("M16").Value = 14.285714285714
$ws.Range("N16").Value = -88.742304309586
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 155
$ws.Range("J17").Value = 126
$ws.Range("K17").Value = 23.015873015873
$ws.Range("L17").Value = 47.619047619047
$ws.Range("M17").Value = 38.392857142857
$ws.Range("N17").Value = -62.918660287081
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 186
$ws.Range("J18").Value = 169
$ws.Range("K18").Value = 10.059171597633
$ws.Range("L18").Value = -16.591928251121
$ws.Range("M18").Value = 13.414634146341
$ws.Range("N18").Value = -89.614740368509
$ws.Range("C19").Value = 53
$ws.Range("D19").Value = 42
$ws.Range("E19").Value = 26.190476190476
$ws.Range("F19").Value = 188
$ws.Range("G19").Value = 137
$ws.Range("H19").Value = 37.226277372262
$ws.Range("I19").Value = 1631
$ws.Range("J19").Value = 974
$ws.Range("K19").Value = 67.453798767967
$ws.Range("L19").Value = 78.056768558952
$ws.Range("M19").Value = 16.416845110635
$ws.Range("N19").Value = -73.866367569299
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 37.5
$ws.Range("J20").Value = 55
$ws.Range("K20").Value = 87.272727272727
$ws.Range("L20").Value = 80.701754385964
$ws.Range("M20").Value = 157.5
$ws.Range("N20").Value = -75.060532687651
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = 22
$ws.Range("F21").Value = 239
$ws.Range("G21").Value = 190
$ws.Range("H21").Value = 25.789473684210
$ws.Range("I21").Value = 2221
$ws.Range("J21").Value = 1458
$ws.Range("K21").Value = 52.331961591220
$ws.Range("L21").Value = 55.859649122807
$ws.Range("M21").Value = 20.379403794037
$ws.Range("N21").Value = -77.860845295055
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 56
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = 36.585365853658
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = 14.285714285714
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 260
$ws.Range("G24").Value = 207
$ws.Range("H24").Value = 25.603864734299
$ws.Range("I24").Value = 2249
$ws.Range("J24").Value = 1606
$ws.Range("K24").Value = 40.037359900373
$ws.Range("L24").Value = 91.241496598639
$ws.Range("M24").Value = 40.124610591900
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 481
$ws.Range("J25").Value = 374
$ws.Range("K25").Value = 28.609625668449
$ws.Range("L25").Value = 81.509433962264
$ws.Range("M25").Value = 27.925531914893
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 66.666666666666
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 12
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 82
$ws.Range("J27").Value = 74
$ws.Range("K27").Value = 10.810810810810
$ws.Range("L27").Value = 57.692307692307
$ws.Range("H30").Value = -100

$excel.CutCopyMode = $false